# "handler for pop up in chrome"
# - Sheet1 ("Places") is trimmed down to a single new entry ("Duren"),
#   dropping the old Madrid/bilbao/valencia rows and the now-empty
#   coordinate cell next to "Duren".
# - Sheet1 becomes the active/selected tab (was Sheet2), with C2 selected.
# - Sheet2 stops being the active tab; its selection moves to F22.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2 is the currently active sheet; update its selection first,
#     while it is still active (selecting a range on the active sheet does
#     not change which tab is active).
[void]$ws2.Range("F22").Select()

# --- Rework Sheet1's data: remove rows 3 & 4 entirely, clear the
#     coordinate value out of row 2, and replace the place name with the
#     new "Duren" entry.
$ws1.Rows("3:4").Delete()
$ws1.Range("B2").ClearContents()
$ws1.Range("A2").Value = "Duren"

# --- Make Sheet1 the active tab and select C2 there (last, so it "wins"
#     the active-tab slot in the saved workbook).
[void]$ws1.Select()
[void]$ws1.Range("C2").Select()
